# Penalty Reward System (unfinished) — trims the oldest and newest data
# points from both the weekly and monthly PO-quantity trend sheets.
#
# Sheet "Weekly Quantity": drop the first 2 data rows and the last 2 data
# rows (26 data rows -> 22 data rows; dimension A1:B27 -> A1:B23).
# Sheet "Monthly Trend": drop the first 1 data row and the last 1 data row
# (9 data rows -> 7 data rows; dimension A1:B10 -> A1:B8).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows.Item(27).Delete()
$ws1.Rows.Item(26).Delete()
$ws1.Rows.Item(2).Delete()
$ws1.Rows.Item(2).Delete()

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Rows.Item(10).Delete()
$ws2.Rows.Item(2).Delete()
